# Updated cryptos list on Sat Jul  8 13:25:14 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto list, and fixes the PaxDollar / TheSandbox rows which had swapped
# their relative order (rows 44 and 45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => updated Price (column D) value. Many of these strings look like
# numbers (and some even contain extra "thousands" dots, e.g. "30.215.14"),
# so force the cells to Text format first to keep the exact original
# textual representation (avoids Excel silently turning "1.000" into 1,
# "97.50" into 97.5, "0.000007371" into scientific notation, etc).
$priceUpdates = [ordered]@{
    2  = "30.215.14"
    3  = "1.859.07"
    4  = "1.000"
    5  = "236.47"
    6  = "1.000"
    8  = "0.2864"
    9  = "0.06538"
    10 = "21.84"
    11 = "0.07930"
    12 = "97.50"
    13 = "1.865.53"
    14 = "5.166"
    15 = "0.6806"
    16 = "266.34"
    17 = "30.207.30"
    18 = "13.78"
    19 = "1.000"
    20 = "0.000007371"
    21 = "2.113.21"
    22 = "5.327"
    23 = "0.9998"
    24 = "6.198"
    25 = "167.25"
    26 = "9.218"
    27 = "18.89"
    28 = "1.966"
    29 = "1.385"
    30 = "0.09862"
    31 = "4.390"
    32 = "1.477"
    33 = "4.067"
    34 = "0.04717"
    35 = "1.129"
    36 = "0.7030"
    37 = "2.707"
    38 = "0.01882"
    39 = "2.631"
    40 = "6.268"
    41 = "74.11"
    42 = "1.943"
    43 = "0.8485"
    44 = "0.4169"
    45 = "0.9995"
    46 = "103.46"
    47 = "966.64"
    48 = "7.171"
    49 = "9.258"
    50 = "34.14"
    51 = "0.05644"
}

# Row => updated Volume(1h) (column E) value (kept padded exactly like the
# original cells: two leading and two trailing spaces).
$volumeUpdates = [ordered]@{
    2  = "  +0.01%  "
    3  = "  -0.27%  "
    4  = "  +0.05%  "
    5  = "  +0.74%  "
    6  = "  +0.09%  "
    7  = "  +0.13%  "
    8  = "  +1.16%  "
    9  = "  +0.00%  "
    10 = "  +4.69%  "
    11 = "  +0.81%  "
    12 = "  +1.24%  "
    13 = "  +0.03%  "
    14 = "  +0.80%  "
    15 = "  +1.28%  "
    16 = "  -4.53%  "
    17 = "  -0.04%  "
    18 = "  +9.22%  "
    20 = "  +1.80%  "
    21 = "  +0.52%  "
    22 = "  -2.55%  "
    23 = "  -0.02%  "
    24 = "  +0.54%  "
    25 = "  +1.47%  "
    26 = "  -0.74%  "
    27 = "  -0.89%  "
    28 = "  +3.60%  "
    29 = "  +2.41%  "
    30 = "  +2.65%  "
    31 = "  -0.06%  "
    32 = "  +0.60%  "
    33 = "  -0.66%  "
    34 = "  +0.55%  "
    35 = "  +2.85%  "
    36 = "  +0.45%  "
    37 = "  -0.26%  "
    38 = "  +1.51%  "
    39 = "  +4.04%  "
    40 = "  -1.61%  "
    41 = "  +1.31%  "
    42 = "  +0.59%  "
    43 = "  +0.43%  "
    44 = "  +0.09%  "
    45 = "  -0.01%  "
    46 = "  -0.22%  "
    47 = "  +5.02%  "
    48 = "  +0.21%  "
    49 = "  +0.02%  "
    50 = "  +0.60%  "
    51 = "  +0.40%  "
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}

# Rows 44 and 45: PaxDollar and TheSandbox swapped places in the ranking.
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
